$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'name'
$ws.Range("B1").Value = 'email'
$ws.Range("C1").Value = 'phone'
$ws.Range("D1").Value = 'gstin'
$ws.Range("E1").Value = 'billing_address'
$ws.Range("F1").Value = 'shipping_address'
$ws.Range("G1").Value = 'notes'
$ws.Range("H1").Value = 'id'
$ws.Range("A2").Value = 'Customer QA'
$ws.Range("B2").Value = 'qa@example.com'
$ws.Range("C2").Value = '''+911234567890'
$ws.Range("D2").Value = 'GST9901QA'
$ws.Range("E2").Value = 'QATown 1'
$ws.Range("F2").Value = 'QATown 2'
$ws.Range("G2").Value = 'Fake customer (test)'
$ws.Range("H2").Value = 'af892bfb-eb9d-40aa-b377-20bb463398bc'
$ws.Range("A3").Value = 'Hassan mansuri'
$ws.Range("B3").Value = 'hassanmansuri570@gmail.com'
$ws.Range("C3").Value = '''+919322909257'
$ws.Range("D3").Value = '''34353'
$ws.Range("E3").Value = 'KPKD'
$ws.Range("F3").Value = 'NAGPUR'
$ws.Range("G3").Value = 'fdghrha'
$ws.Range("H3").Value = '57fe89c5-a399-4dd1-9830-f513fc466f73'
$ws.Range("A4").Value = 'gogo'
$ws.Range("B4").Value = 'hassanmansuri570@gmail.com'
$ws.Range("C4").Value = '+919322909257ee'
$ws.Range("D4").Value = ''
$ws.Range("E4").Value = 'KPKDefe'
$ws.Range("F4").Value = 'NAGPURef'
$ws.Range("G4").Value = 'fdbfbfdb e gh'
$ws.Range("H4").Value = '6a4f98a3-e2fa-44da-8b37-868632dfbee8'
$ws.Range("A5").Value = 'Customer QA'
$ws.Range("B5").Value = 'qa@example.com'
$ws.Range("C5").Value = '''+911234567890'
$ws.Range("D5").Value = 'GST9901QA'
$ws.Range("E5").Value = 'QATown 1'
$ws.Range("F5").Value = 'QATown 2'
$ws.Range("G5").Value = 'Fake customer (test)'
$ws.Range("H5").Value = '7dd73460-a560-4874-886f-78e863a66d49'
